# Daily attendance processing - normalize "Recorded By" (column G) entries
# so that the "System" marker is listed first among the recorder names.
#
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#          "backup@backdoor.com, system, System" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Exact-match replacements observed for the "Recorded By" column (column G).
$recordedByFixes = @{
    "dnasr281@gmail.com, System"            = "System, dnasr281@gmail.com"
    "admin@admin.com, System"               = "System, admin@admin.com"
    "backup@backdoor.com, system, System"   = "system, backup@backdoor.com, System"
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2

    if ($current -ne $null -and $recordedByFixes.ContainsKey($current)) {
        $cell.Value2 = $recordedByFixes[$current]
    }
}
